$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 78.09999999999999
$ws.Range("B4").Value = 90.59999999999999
$ws.Range("C4").Value = 33.3
$ws.Range("B5").Value = 94.7
$ws.Range("B6").Value = 87.3
